# Auto-generated edit script updating crypto price/volume columns (D, E)
# to reflect refreshed values from the GitHub Actions data-refresh run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.300.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.21%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.867.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.91%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'243.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.38%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.01%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4729"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.70%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2867"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.48%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.17%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'21.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.85%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.76%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'97.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.51%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.867.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.89%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.7219"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.55%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -2.00%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'280.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.69%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'30.280.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.57%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'13.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.37%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.0000"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.000007459"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.38%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.111.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.89%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.236"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.05%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.258"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.21%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'162.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.23%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.008"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.61%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.90%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.878"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.13%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.09617"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "'  -2.21%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.476"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.89%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.219"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.76%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.55%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.04788"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.61%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.118"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.87%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.6847"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.36%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E38").Value = "'  -1.41%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.831"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.42%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'75.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.62%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.206"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.85%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.934"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -4.80%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.4205"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.31%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.9990"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.05%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.8244"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.19%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'100.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.22%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'9.592"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.11%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'6.951"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.66%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'34.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.10%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.05763"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.02%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'882.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.87%  "
$ws.Range("E51").Style = "Normal"
